$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K1").Value = 0.95690985494037928
$ws.Range("AC1").Value = 0.69316183565096778
$ws.Range("AJ1").Value = 0.61985945480376847
$ws.Range("AY1").Value = 0.98433942637759664
$ws.Range("BA1").Value = 0.63238974420440841
$ws.Range("D2").Value = 0.75260813441963514
$ws.Range("Z2").Value = 0.72072965909210329
$ws.Range("Q3").Value = 0.59246547643864533
$ws.Range("AL3").Value = 0.62811558786603228
$ws.Range("BL3").Value = 0.80258167972251249
$ws.Range("BN3").Value = 0.76345798019658084
$ws.Range("B4").Value = 0.68386445437095877
$ws.Range("P4").Value = 0.6149342684510436
$ws.Range("AA4").Value = 0.67229327639129755
$ws.Range("AG4").Value = 0.95469166287526086
$ws.Range("BC4").Value = 0.91933216460947631
$ws.Range("F5").Value = 0.88326473884509982
$ws.Range("BP5").Value = 0.80009068185504184
$ws.Range("E6").Value = 0.67752855301384818
$ws.Range("J6").Value = 0.60630877821746909
$ws.Range("Z6").Value = 0.51468603331710927
$ws.Range("AP6").Value = 0.92771599813516259
$ws.Range("AS6").Value = 0.8130053690313841
$ws.Range("BE6").Value = 0.62933927662447808
$ws.Range("BP6").Value = 0.94786956920466614
$ws.Range("AE7").Value = 0.76495589797256924
$ws.Range("BH7").Value = 0.52215408291971577
$ws.Range("BM7").Value = 0.86294856075762683
$ws.Range("Q8").Value = 0.80694364384154738
$ws.Range("BE8").Value = 0.52169462464166494
$ws.Range("Q9").Value = 0.9688787403338639
$ws.Range("BJ9").Value = 0.54269960923627791
$ws.Range("F10").Value = 0.68784012434306141
$ws.Range("L10").Value = 0.60740234934435067
$ws.Range("AB10").Value = 0.53845202540423176
$ws.Range("AK10").Value = 0.86159846851612754
$ws.Range("A11").Value = 0.6246009803620407
$ws.Range("Z11").Value = 0.84497474179232823
$ws.Range("AC11").Value = 0.61807592779928
$ws.Range("AK11").Value = 0.63685312678698947
$ws.Range("J12").Value = 0.65746289149034975
$ws.Range("M12").Value = 0.56138611387544524
$ws.Range("BB12").Value = 0.89262345172684721
$ws.Range("BF12").Value = 0.59883816617705787
$ws.Range("L13").Value = 0.55111193157624871
$ws.Range("Q13").Value = 0.72516483132712084
$ws.Range("W13").Value = 0.88649211384038329
$ws.Range("BH13").Value = 0.62988597218356202
$ws.Range("T14").Value = 0.83698554600281883
$ws.Range("BL14").Value = 0.91549816795919281
$ws.Range("P15").Value = 0.86925743030850611
$ws.Range("X15").Value = 0.54416775001662887
$ws.Range("AS15").Value = 0.80463991625808529
$ws.Range("BL15").Value = 0.67637068469460404
$ws.Range("D16").Value = 0.51867196091633949
$ws.Range("O16").Value = 0.91740445972377649
$ws.Range("S16").Value = 0.57070444663071918
$ws.Range("AD16").Value = 0.52648313991274986
$ws.Range("AQ16").Value = 0.88708640535412697
$ws.Range("BB16").Value = 0.5313290602735028
$ws.Range("C17").Value = 0.88412997211609379
$ws.Range("H17").Value = 0.918174679060211
$ws.Range("I17").Value = 0.84773501742615665
$ws.Range("M17").Value = 0.7354661355896237
$ws.Range("AT17").Value = 0.96955834656204953
$ws.Range("BD17").Value = 0.8991713357900335
$ws.Range("BH17").Value = 0.84657367085622282
$ws.Range("U18").Value = 0.59153832991687483
$ws.Range("AR18").Value = 0.82525297895195693
$ws.Range("P19").Value = 0.91670249082597
$ws.Range("AQ19").Value = 0.59324982760698775
$ws.Range("BM19").Value = 0.7414209176214448
$ws.Range("BP19").Value = 0.89229586056404298
$ws.Range("N20").Value = 0.82428264044146615
$ws.Range("V20").Value = 0.83408972762578459
$ws.Range("W20").Value = 0.6088347985218534
$ws.Range("AM20").Value = 0.80093533050590471
$ws.Range("R21").Value = 0.93877393174646895
$ws.Range("W21").Value = 0.5594907407245463
$ws.Range("AL21").Value = 0.61999194761261855
$ws.Range("T22").Value = 0.89927707136243873
$ws.Range("X22").Value = 0.99850965438256356
$ws.Range("AQ22").Value = 0.63580570980232465
$ws.Range("AR22").Value = 0.95993916944724855
$ws.Range("M23").Value = 0.85800996884594705
$ws.Range("T23").Value = 0.64054706970900266
$ws.Range("U23").Value = 0.99734986394956804
$ws.Range("BG23").Value = 0.95908520810961184
$ws.Range("BN23").Value = 0.9227183287470696
$ws.Range("O24").Value = 0.70922902305980062
$ws.Range("V24").Value = 0.72789652578344732
$ws.Range("AA24").Value = 0.62984109754040296
$ws.Range("AE25").Value = 0.79396587518027251
$ws.Range("AM25").Value = 0.84087661991869145
$ws.Range("AO25").Value = 0.72641786645862716
$ws.Range("B26").Value = 0.76034028781834673
$ws.Range("F26").Value = 0.93088419505571562
$ws.Range("K26").Value = 0.82283512426622374
$ws.Range("AC26").Value = 0.63956056317957966
$ws.Range("BJ26").Value = 0.68332753841195237
$ws.Range("D27").Value = 0.86906455607161548
$ws.Range("X27").Value = 0.62116705704821273
$ws.Range("AG27").Value = 0.89860629356920518
$ws.Range("BD27").Value = 0.92326450030134777
$ws.Range("BG27").Value = 0.57983074589084938
$ws.Range("J28").Value = 0.77542884303300719
$ws.Range("AD28").Value = 0.63443968560908881
$ws.Range("AV28").Value = 0.84602676093977425
$ws.Range("A29").Value = 0.60486212088501956
$ws.Range("K29").Value = 0.57555811608005247
$ws.Range("Z29").Value = 0.60579591332768756
$ws.Range("AW29").Value = 0.98916753483097264
$ws.Range("BK29").Value = 0.85551625414070909
$ws.Range("P30").Value = 0.72585122791354495
$ws.Range("AB30").Value = 0.98840462531431528
$ws.Range("BA30").Value = 0.92583067579656131
$ws.Range("G31").Value = 0.67121790144119076
$ws.Range("Y31").Value = 0.78929078187219037
$ws.Range("BD31").Value = 0.76881603849693336
$ws.Range("AG32").Value = 0.80862069463297614
$ws.Range("BA32").Value = 0.91490820171600862
$ws.Range("D33").Value = 0.5042648895152142
$ws.Range("AA33").Value = 0.9205474264493021
$ws.Range("AF33").Value = 0.6225278845751655
$ws.Range("BJ33").Value = 0.65748730327194727
$ws.Range("AU34").Value = 0.7273323702304636
$ws.Range("AV34").Value = 0.97491522085995186
$ws.Range("AX34").Value = 0.63738709448642572
$ws.Range("BJ34").Value = 0.88108319337708685
$ws.Range("BE35").Value = 0.51681678888647087
$ws.Range("BN35").Value = 0.51483581700988934
$ws.Range("A36").Value = 0.63594335259950463
$ws.Range("AK36").Value = 0.63403051218526785
$ws.Range("AZ36").Value = 0.89834601505221834
$ws.Range("J37").Value = 0.881506940997969
$ws.Range("K37").Value = 0.98744276656462193
$ws.Range("AJ37").Value = 0.52116947904701694
$ws.Range("BG37").Value = 0.99500461043170962
$ws.Range("BN37").Value = 0.584657110365707
$ws.Range("C38").Value = 0.52158032908315355
$ws.Range("U38").Value = 0.80250584074073883
$ws.Range("AX38").Value = 0.71713183421124427
$ws.Range("BB38").Value = 0.68352729528287282
$ws.Range("T39").Value = 0.67965653794501413
$ws.Range("Y39").Value = 0.55426936376596325
$ws.Range("AU39").Value = 0.71096343048941446
$ws.Range("BC39").Value = 0.65342912384962881
$ws.Range("BI39").Value = 0.85543132098001728
$ws.Range("AO40").Value = 0.55780936821625948
$ws.Range("AW40").Value = 0.89818213089559662
$ws.Range("BI40").Value = 0.68975880823185998
$ws.Range("Y41").Value = 0.94899645336515182
$ws.Range("AN41").Value = 0.86026205006888645
$ws.Range("AQ41").Value = 0.86954398621813223
$ws.Range("AX41").Value = 0.92371305388975877
$ws.Range("F42").Value = 0.62461609580740518
$ws.Range("AR42").Value = 0.70106933468869181
$ws.Range("P43").Value = 0.96442455386966741
$ws.Range("S43").Value = 0.92312751556160033
$ws.Range("V43").Value = 0.97632734518115316
$ws.Range("AO43").Value = 0.85576223335153534
$ws.Range("AR43").Value = 0.99449564197577711
$ws.Range("R44").Value = 0.57008120583449617
$ws.Range("V44").Value = 0.65296241830423885
$ws.Range("AP44").Value = 0.66210874358330374
$ws.Range("AQ44").Value = 0.59759886785624294
$ws.Range("AS44").Value = 0.75552372794171929
$ws.Range("F45").Value = 0.56011807254675239
$ws.Range("O45").Value = 0.6258997298861827
$ws.Range("AR45").Value = 0.82625584525566687
$ws.Range("Q46").Value = 0.86531601242195566
$ws.Range("AU46").Value = 0.50521581413460903
$ws.Range("AV46").Value = 0.99260261616495926
$ws.Range("AY46").Value = 0.91434726012884426
$ws.Range("BF46").Value = 0.66631791309792809
$ws.Range("AH47").Value = 0.87464447294438918
$ws.Range("AM47").Value = 0.56738974083162075
$ws.Range("AT47").Value = 0.81148141624614034
$ws.Range("BB47").Value = 0.76617013000726875
$ws.Range("AB48").Value = 0.56885797212975875
$ws.Range("AH48").Value = 0.85946098675496763
$ws.Range("AT48").Value = 0.73358589951130382
$ws.Range("AZ48").Value = 0.55139416333282854
$ws.Range("BL48").Value = 0.65317219771886503
$ws.Range("BO48").Value = 0.71785746411106055
$ws.Range("AC49").Value = 0.72223394444141209
$ws.Range("AN49").Value = 0.68132348476578031
$ws.Range("BB49").Value = 0.66384673229288005
$ws.Range("AH50").Value = 0.52212454163345123
$ws.Range("AL50").Value = 0.95526535169703819
$ws.Range("AO50").Value = 0.5052888384641403
$ws.Range("AZ50").Value = 0.72098263271285168
$ws.Range("BC50").Value = 0.61325298889363633
$ws.Range("A51").Value = 0.6079783049655767
$ws.Range("AT51").Value = 0.61713398614271076
$ws.Range("BA51").Value = 0.59352738972667041
$ws.Range("AJ52").Value = 0.85377642925819819
$ws.Range("AV52").Value = 0.8648616443933097
$ws.Range("AX52").Value = 0.82167331734015026
$ws.Range("BM52").Value = 0.77366998178157398
$ws.Range("A53").Value = 0.90284238102449921
$ws.Range("AD53").Value = 0.95656455500286874
$ws.Range("AF53").Value = 0.52335706661982773
$ws.Range("AY53").Value = 0.59954244769448795
$ws.Range("BJ53").Value = 0.82385374146262025
$ws.Range("L54").Value = 0.79888285118565006
$ws.Range("P54").Value = 0.61219552436300995
$ws.Range("AL54").Value = 0.99631264498601591
$ws.Range("AU54").Value = 0.55485825624282581
$ws.Range("AW54").Value = 0.89362242487768306
$ws.Range("BL54").Value = 0.95940569668025488
$ws.Range("D55").Value = 0.93426710274319102
$ws.Range("AM55").Value = 0.97245828478223217
$ws.Range("AX55").Value = 0.60802828009339283
$ws.Range("BE55").Value = 0.64636133193809753
$ws.Range("Q56").Value = 0.57369478232134152
$ws.Range("AA56").Value = 0.86150639483052061
$ws.Range("AE56").Value = 0.53701660720333755
$ws.Range("F57").Value = 0.958380338915507
$ws.Range("H57").Value = 0.805096372645695
$ws.Range("AI57").Value = 0.74531744757373519
$ws.Range("BC57").Value = 0.79477306118386548
$ws.Range("BH57").Value = 0.51910060265172353
$ws.Range("BN57").Value = 0.58542566961273756
$ws.Range("L58").Value = 0.70731775012144338
$ws.Range("AT58").Value = 0.73916412738760273
$ws.Range("BO58").Value = 0.64958350862731362
$ws.Range("W59").Value = 0.61402259154695848
$ws.Range("AA59").Value = 0.85401241498434033
$ws.Range("AK59").Value = 0.98209845320894384
$ws.Range("G60").Value = 0.55855215746132225
$ws.Range("M60").Value = 0.51599580545836665
$ws.Range("Q60").Value = 0.83567764027839186
$ws.Range("BE60").Value = 0.62760639110425886
$ws.Range("BM60").Value = 0.95678091457435499
$ws.Range("AM61").Value = 0.9765504656458418
$ws.Range("AN61").Value = 0.72204500258382565
$ws.Range("BJ61").Value = 0.6495536603738894
$ws.Range("I62").Value = 0.98562167002482015
$ws.Range("Z62").Value = 0.84086917317112397
$ws.Range("AG62").Value = 0.84451791304443358
$ws.Range("AH62").Value = 0.52474058180907568
$ws.Range("BA62").Value = 0.84151513735749284
$ws.Range("BI62").Value = 0.86247258627892576
$ws.Range("BL62").Value = 0.79690261211169511
$ws.Range("BP62").Value = 0.60071912215119116
$ws.Range("AC63").Value = 0.5657075483364542
$ws.Range("BL63").Value = 0.95452493736322075
$ws.Range("C64").Value = 0.67128569060069077
$ws.Range("N64").Value = 0.8009831671377281
$ws.Range("O64").Value = 0.70764503488811448
$ws.Range("AV64").Value = 0.6628546457426896
$ws.Range("BB64").Value = 0.57128728776580728
$ws.Range("BJ64").Value = 0.68229120082675077
$ws.Range("BK64").Value = 0.64673813746079534
$ws.Range("BM64").Value = 0.7174543367874705
$ws.Range("G65").Value = 0.93295248560825716
$ws.Range("S65").Value = 0.50700371382214071
$ws.Range("AZ65").Value = 0.50302001257475382
$ws.Range("BH65").Value = 0.6093800056186004
$ws.Range("BL65").Value = 0.60810663502515383
$ws.Range("C66").Value = 0.95621432154440011
$ws.Range("W66").Value = 0.62577664986808967
$ws.Range("AI66").Value = 0.93734824790173199
$ws.Range("AK66").Value = 0.66483727359160971
$ws.Range("BE66").Value = 0.97643321640431768
$ws.Range("AV67").Value = 0.54076455891655617
$ws.Range("BF67").Value = 0.87672626743772686
$ws.Range("E68").Value = 0.56653173963400416
$ws.Range("F68").Value = 0.97346520613461274
$ws.Range("S68").Value = 0.81188911454575297
$ws.Range("BJ68").Value = 0.58856798639837082

# Column width changes
for ($i = 1; $i -le 68; $i++) {
    if ($i -eq 5) {
        $ws.Columns.Item($i).ColumnWidth = 10.833333333333332
    } else {
        $ws.Columns.Item($i).ColumnWidth = 11.833333333333332
    }
}